$wb = $excel.ActiveWorkbook

# Source sheet (PCA_100_Polarity) has the header/format layout we want to
# replicate on the target sheet (PCA-300-Polarity).
$wsSrc = $wb.Worksheets.Item(4)
$ws = $wb.Worksheets.Item(5)

# Copy the cell formatting (borders, bold, wrap text) from the analogous
# table on sheet 4 so the new header/data cells reuse the same styles.
$wsSrc.Range("A1:M4").Copy()
$ws.Range("A1:M4").PasteSpecial(-4122)

# Header row (row 1)
$ws.Range("A1").Value = "ngrams"
$ws.Range("B1").Value = "logreg_train"
$ws.Range("C1").Value = "logreg_test"
$ws.Range("D1").Value = " nb_train"
$ws.Range("E1").Value = " nb_test"
$ws.Range("F1").Value = " svm_train"
$ws.Range("G1").Value = " svm_test"
$ws.Range("H1").Value = " dt_train"
$ws.Range("I1").Value = " dt_test"
$ws.Range("J1").Value = " adaboost_train"
$ws.Range("K1").Value = " adaboost_test"
$ws.Range("L1").Value = " rf_train"
$ws.Range("M1").Value = " rf_test"

# Row 2 (1gram)
$ws.Range("A2").Value = "1gram"
$ws.Range("B2").Value = 0.76173900000000005
$ws.Range("C2").Value = 0.75984300000000005
$ws.Range("D2").Value = 0.58704500000000004
$ws.Range("E2").Value = 0.58754200000000001
$ws.Range("F2").Value = 0.76102000000000003
$ws.Range("G2").Value = 0.76074699999999995
$ws.Range("H2").Value = 0.634405
$ws.Range("I2").Value = 0.63171299999999997
$ws.Range("J2").Value = 0.776868
$ws.Range("K2").Value = 0.75904300000000002
$ws.Range("L2").Value = 0.67621200000000004
$ws.Range("M2").Value = 0.67383099999999996

# Row 3 (2gram)
$ws.Range("A3").Value = "2gram"
$ws.Range("B3").Value = 0.60944399999999999
$ws.Range("C3").Value = 0.60322799999999999
$ws.Range("D3").Value = 0.55218299999999998
$ws.Range("E3").Value = 0.55161400000000005
$ws.Range("F3").Value = 0.60821499999999995
$ws.Range("G3").Value = 0.60246200000000005
$ws.Range("H3").Value = 0.57006100000000004
$ws.Range("I3").Value = 0.56510899999999997
$ws.Range("J3").Value = 0.63682799999999995
$ws.Range("K3").Value = 0.60987100000000005
$ws.Range("L3").Value = 0.589202
$ws.Range("M3").Value = 0.582986

# Row 4 (3gram)
$ws.Range("A4").Value = "3gram"
$ws.Range("B4").Value = 0.53217199999999998
$ws.Range("C4").Value = 0.52288500000000004
$ws.Range("D4").Value = 0.51551199999999997
$ws.Range("E4").Value = 0.50866
$ws.Range("F4").Value = 0.53228799999999998
$ws.Range("G4").Value = 0.52316399999999996
$ws.Range("H4").Value = 0.51734400000000003
$ws.Range("I4").Value = 0.51196399999999997
$ws.Range("J4").Value = 0.54363899999999998
$ws.Range("K4").Value = 0.52660700000000005
$ws.Range("L4").Value = 0.53046800000000005
$ws.Range("M4").Value = 0.52097199999999999

# Row heights
$ws.Rows.Item(1).RowHeight = 46
$ws.Rows.Item(2).RowHeight = 17
$ws.Rows.Item(3).RowHeight = 17
$ws.Rows.Item(4).RowHeight = 17

# Match the saved selection state (J2:K2 active)
$null = $ws.Range("J2:K2").Select()
